## feat: add 2022-Q1 data
#
# Adds a new "2022-Q1" holdings sheet (same layout as the other quarterly
# sheets) positioned between "2021-Q1" and "总计", and records the new
# quarter as the newest row of the "总计" (totals) rollup sheet.

$wb = $excel.ActiveWorkbook

$q1sheet = $wb.Worksheets.Item("2021-Q1")
$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. Create the "2022-Q1" sheet -----------------------------------
# Duplicate "总计" (rather than Worksheets.Add()) so the new sheet starts
# out sharing the existing header/row formatting (bold + thin border) -
# pasting formats onto a brand-new blank sheet later in the same script
# doesn't stick, but a sheet clone keeps its styles from the get-go.
$totalSheet.Copy($null, $q1sheet)
$newSheet = $wb.Worksheets.Item("总计 (2)")
$newSheet.Name = "2022-Q1"

# "总计" has only 4 columns (A:D); this layout needs 8 (A:H), so extend the
# header's formatting (bold, centered, thin border) across E1:H1.
$newSheet.Range("D1").Copy()
$newSheet.Range("E1:H1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Row 2 holds the single 2022-Q1 fund holding. Only A2 keeps the bold
# bordered style inherited from the clone - the data cells (B2:H2) are
# unstyled, so strip whatever they inherited from the copied sheet.
$newSheet.Range("B2:H2").Style = "Normal"

$newSheet.Range("A2").Value = 0

# Leading "'" forces text storage (matches the source data, which keeps
# these as strings rather than numbers); Style="Normal" afterwards clears
# the "Text" number-format Excel stamps on when you type a quoted number,
# so the cell ends up with no explicit style - same as the target.
$newSheet.Range("B2").Value = "'310368"
$newSheet.Range("B2").Style = "Normal"

$newSheet.Range("C2").Value = "申万菱信竞争优势混合"

$newSheet.Range("D2").Value = "'0.83"
$newSheet.Range("D2").Style = "Normal"

$newSheet.Range("E2").Value = "'91.22"
$newSheet.Range("E2").Style = "Normal"

$newSheet.Range("F2").Value = "'4.53"
$newSheet.Range("F2").Style = "Normal"

$newSheet.Range("G2").Value = "'0.0376"
$newSheet.Range("G2").Style = "Normal"

$newSheet.Range("H2").Value = 2

# --- 2. Record the new quarter at the top of "总计" ---------------------
# Re-fetch "总计" by name: worksheet handles are positional, and the
# sheet-copy/insert above shifted the collection, so the stale $totalSheet
# reference would now resolve to the wrong tab.
$totalSheet = $wb.Worksheets.Item("总计")

$totalSheet.Range("A2").EntireRow.Insert()

# New row 2 inherits the header's bold/border style on insert; A2 should
# instead match A3 (plain bold+border, no background), and B2:D2 should
# carry no style at all - same pattern as the "2022-Q1" sheet above.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").Style = "Normal"

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.04

$totalSheet.Range("A3").Value = 1

# Leave the original sheet ("2021-Q1") selected, matching the unmodified
# <bookViews>/tabSelected state from before this edit.
$q1sheet.Activate()
